# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the refreshed scrape output.
$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1795
$ws1.Range("F4").Value = 461
$ws1.Range("F8").Value = 342
$ws1.Range("F9").Value = 1746
$ws1.Range("F10").Value = 373
$ws1.Range("F11").Value = 1430
$ws1.Range("F13").Value = 341
$ws1.Range("F14").Value = 686
$ws1.Range("F15").Value = 12851
$ws1.Range("F16").Value = 12835
$ws1.Range("F18").Value = 747
$ws1.Range("F20").Value = 520
$ws1.Range("F22").Value = 576
$ws1.Range("F23").Value = 2018
$ws1.Range("F25").Value = 11
$ws1.Range("F27").Value = 68
$ws1.Range("F29").Value = 681

# "演出" sheet (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 128
$ws2.Range("F10").Value = 80

# "本地生活" sheet (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 168

# "全部类型" sheet (sheet4) - combined view of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 168
$ws4.Range("F5").Value = 1795
$ws4.Range("F6").Value = 461
$ws4.Range("F13").Value = 342
$ws4.Range("F14").Value = 1746
$ws4.Range("F15").Value = 373
$ws4.Range("F16").Value = 1430
$ws4.Range("F18").Value = 341
$ws4.Range("F20").Value = 686
$ws4.Range("F21").Value = 12851
$ws4.Range("F22").Value = 12835
$ws4.Range("F24").Value = 747
$ws4.Range("F26").Value = 520
$ws4.Range("F28").Value = 576
$ws4.Range("F31").Value = 2018
$ws4.Range("F33").Value = 11
$ws4.Range("F34").Value = 128
$ws4.Range("F37").Value = 68
$ws4.Range("F39").Value = 681
$ws4.Range("F40").Value = 80
